# Update latest output (run 191)
$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet: updated cost / unit cost totals ---
$wsSchedule.Range("E2").Value = 849.5128297500002
$wsSchedule.Range("F2").Value = 14.0461777405754

# --- Detailed sheet: updated price values (and a few Type changes) ---
$wsDetailed.Range("B21").Value = 0.51

$wsDetailed.Range("B22").Value = 0.02957

$wsDetailed.Range("B23").Value = 0.7
$wsDetailed.Range("C23").Value = "historical"

$wsDetailed.Range("B24").Value = 5.26268
$wsDetailed.Range("C24").Value = "historical"

$wsDetailed.Range("B25").Value = 0.01783
$wsDetailed.Range("C25").Value = "historical"

$wsDetailed.Range("B26").Value = 0
$wsDetailed.Range("C26").Value = "historical"

$wsDetailed.Range("B27").Value = 0

$wsDetailed.Range("B28").Value = -7.2

$wsDetailed.Range("B29").Value = -12.01

$wsDetailed.Range("B30").Value = -11.01

$wsDetailed.Range("B31").Value = -10

$wsDetailed.Range("B32").Value = -8.524710000000001

$wsDetailed.Range("B33").Value = -6.83333

$wsDetailed.Range("B34").Value = -4.89981

$wsDetailed.Range("B35").Value = 0.51

$wsDetailed.Range("B36").Value = 9.230259999999999

$wsDetailed.Range("B37").Value = 22.73226

$wsDetailed.Range("B38").Value = 46.10125

$wsDetailed.Range("B39").Value = 55.64524

$wsDetailed.Range("B41").Value = 60.16225

$wsDetailed.Range("B42").Value = 64.89

$wsDetailed.Range("B43").Value = 61.2959

$wsDetailed.Range("B44").Value = 63.23165

$wsDetailed.Range("B45").Value = 62.33147

$wsDetailed.Range("B49").Value = 57.06
